$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Accion" rows (A2:B10) get reshuffled: the first data row ("No ocurre
# nada", id 0) moves to the bottom of the list, every other row shifts up
# one position, and the IDs are renumbered 1..9 instead of 0..8.

$labels = @(
    "Quita X puntos de vida (primero ataca al escudo)",
    "Ignorar escudos (ataca directamente al casco)",
    "Destruir escudo",
    "Enviar a enfermeria",
    "Tirar dado amenaza",
    "Amenazas externas recuperan 1 nivel ",
    "Devolver dado bloqueado",
    "No asignar",
    "No ocurre nada"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $labels[$i]
}

$ws.Range("D8").Select()
